# Apply the edits described by the diff: update the date and each
# three-digit x one-digit multiplication fact in the table cells.
$d = $word.ActiveDocument

$d.Content.Find.Execute('2025-02-07 Friday', $true, $false, $false, $false, $false, $true, 1, $false, '2025-02-08 Saturday', 2) | Out-Null
$d.Content.Find.Execute('284×2=568', $true, $false, $false, $false, $false, $true, 1, $false, '128×3=384', 2) | Out-Null
$d.Content.Find.Execute('577×2=1154', $true, $false, $false, $false, $false, $true, 1, $false, '452×7=3164', 2) | Out-Null
$d.Content.Find.Execute('179×2=358', $true, $false, $false, $false, $false, $true, 1, $false, '516×7=3612', 2) | Out-Null
$d.Content.Find.Execute('448×2=896', $true, $false, $false, $false, $false, $true, 1, $false, '646×8=5168', 2) | Out-Null
$d.Content.Find.Execute('216×3=648', $true, $false, $false, $false, $false, $true, 1, $false, '527×3=1581', 2) | Out-Null
$d.Content.Find.Execute('910×9=8190', $true, $false, $false, $false, $false, $true, 1, $false, '337×8=2696', 2) | Out-Null
$d.Content.Find.Execute('192×8=1536', $true, $false, $false, $false, $false, $true, 1, $false, '813×4=3252', 2) | Out-Null
$d.Content.Find.Execute('484×6=2904', $true, $false, $false, $false, $false, $true, 1, $false, '139×5=695', 2) | Out-Null
$d.Content.Find.Execute('443×4=1772', $true, $false, $false, $false, $false, $true, 1, $false, '845×4=3380', 2) | Out-Null
$d.Content.Find.Execute('424×2=848', $true, $false, $false, $false, $false, $true, 1, $false, '440×9=3960', 2) | Out-Null
$d.Content.Find.Execute('862×3=2586', $true, $false, $false, $false, $false, $true, 1, $false, '488×9=4392', 2) | Out-Null
$d.Content.Find.Execute('274×3=822', $true, $false, $false, $false, $false, $true, 1, $false, '342×3=1026', 2) | Out-Null
$d.Content.Find.Execute('785×4=3140', $true, $false, $false, $false, $false, $true, 1, $false, '121×5=605', 2) | Out-Null
$d.Content.Find.Execute('268×5=1340', $true, $false, $false, $false, $false, $true, 1, $false, '867×2=1734', 2) | Out-Null
$d.Content.Find.Execute('102×3=306', $true, $false, $false, $false, $false, $true, 1, $false, '548×2=1096', 2) | Out-Null
$d.Content.Find.Execute('895×7=6265', $true, $false, $false, $false, $false, $true, 1, $false, '867×4=3468', 2) | Out-Null
$d.Content.Find.Execute('913×2=1826', $true, $false, $false, $false, $false, $true, 1, $false, '409×8=3272', 2) | Out-Null
$d.Content.Find.Execute('129×4=516', $true, $false, $false, $false, $false, $true, 1, $false, '186×7=1302', 2) | Out-Null
$d.Content.Find.Execute('941×5=4705', $true, $false, $false, $false, $false, $true, 1, $false, '747×3=2241', 2) | Out-Null
$d.Content.Find.Execute('705×5=3525', $true, $false, $false, $false, $false, $true, 1, $false, '694×5=3470', 2) | Out-Null
$d.Content.Find.Execute('390×6=2340', $true, $false, $false, $false, $false, $true, 1, $false, '658×7=4606', 2) | Out-Null
$d.Content.Find.Execute('375×2=750', $true, $false, $false, $false, $false, $true, 1, $false, '873×2=1746', 2) | Out-Null
$d.Content.Find.Execute('520×2=1040', $true, $false, $false, $false, $false, $true, 1, $false, '745×6=4470', 2) | Out-Null
$d.Content.Find.Execute('479×6=2874', $true, $false, $false, $false, $false, $true, 1, $false, '942×8=7536', 2) | Out-Null
$d.Content.Find.Execute('331×3=993', $true, $false, $false, $false, $false, $true, 1, $false, '303×2=606', 2) | Out-Null
